$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. resolution_z_unit list sheet: add "mm" as a new first option,
#    keeping existing nm / um, reordered as mm, um, nm.
# ------------------------------------------------------------------
$zlist = $wb.Worksheets.Item("resolution_z_unit list")
$zlist.Range("A3").Value2 = "nm"
$zlist.Range("A2").Value2 = "um"
$zlist.Range("A1").Value2 = "mm"

# ------------------------------------------------------------------
# 2. Main sheet: update comment text and the resolution_z_unit
#    column's data validation to match the new list.
# ------------------------------------------------------------------
$main = $wb.Worksheets.Item("Export as TSV")

$cmt = $main.Range("V1").Comment
[void]$cmt.Text("The unit of incremental distance between image slices.")

$val = $main.Range("V2:V1048576").Validation
$val.Modify(3, 1, 1, "'resolution_z_unit list'!`$A`$1:`$A`$3")
$val.ErrorTitle = "Value must come from list"
$val.ErrorMessage = "Value must be one of: mm / um / nm."
